# Add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" detail sheet right after "总计" (pushing the
#    existing 2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q3 sheets one slot over).
# 2. Insert a matching summary row ("2022-Q4", 2, 0.03) at the top of the
#    "总计" roll-up sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" roll-up sheet: insert a new row right under the header row
#    holding the 2022-Q4 summary figures, shifting the existing quarters
#    down by one row.
# ------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

# The row-insert carries the header row's bold/bordered formatting down
# onto B2:D2 - strip that back off so the new row matches the plain
# (unstyled) look of the other data rows.
$totals.Range($totals.Cells.Item(2,2), $totals.Cells.Item(2,4)).ClearFormats()

# Re-apply the index-column style (bold/centered/bordered, same as the
# style already used on A3:A6) to the new A2 cell.
$totals.Cells.Item(3,1).Copy()
$totals.Cells.Item(2,1).PasteSpecial(-4122)

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q4"
$totals.Cells.Item(2,3).Value = 2
$totals.Cells.Item(2,4).Value = 0.03

# Renumber the index column (col A) of the rows that got pushed down -
# it's a 0-based row counter, so each shifted row's index goes up by 1.
$totals.Cells.Item(3,1).Value = 1
$totals.Cells.Item(4,1).Value = 2
$totals.Cells.Item(5,1).Value = 3
$totals.Cells.Item(6,1).Value = 4

# ------------------------------------------------------------------
# 2) New "2022-Q4" detail sheet, positioned right after "总计".
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($totals.Next)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

# Force the data block to be stored as text (matches the source data,
# which keeps e.g. "001942" and "90.30" as literal strings rather than
# numbers), then clear the formatting residue afterwards so the cells
# end up unstyled just like in the sibling quarterly sheets.
$dataRange = $q4.Range($q4.Cells.Item(2,2), $q4.Cells.Item(3,7))
$dataRange.NumberFormat = "@"

for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "001942"
$q4.Cells.Item(2,3).Value = "前海开源沪港深汇鑫灵活配置混合A"
$q4.Cells.Item(2,4).Value = "0.31"
$q4.Cells.Item(2,5).Value = "90.30"
$q4.Cells.Item(2,6).Value = "4.54"
$q4.Cells.Item(2,7).Value = "0.0141"
$q4.Cells.Item(2,8).Value = 8

$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "001943"
$q4.Cells.Item(3,3).Value = "前海开源沪港深汇鑫灵活配置混合C"
$q4.Cells.Item(3,4).Value = "0.27"
$q4.Cells.Item(3,5).Value = "90.30"
$q4.Cells.Item(3,6).Value = "4.54"
$q4.Cells.Item(3,7).Value = "0.0123"
$q4.Cells.Item(3,8).Value = 8

$dataRange.ClearFormats()

# Re-apply the header-row style (bold/centered/bordered "s=2" style used
# by every sheet's header row, cols B:H) and the index-column style
# (col A, rows 2:3) to match the other quarterly sheets.
$totals.Cells.Item(1,2).Copy()
$q4.Range($q4.Cells.Item(1,2), $q4.Cells.Item(1,8)).PasteSpecial(-4122)
$totals.Cells.Item(2,1).Copy()
$q4.Range($q4.Cells.Item(2,1), $q4.Cells.Item(3,1)).PasteSpecial(-4122)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(3,1).Value = 1

$totals.Select()
$totals.Range("A1").Select()
